# Applies the "Updated with final part quantities and costs" commit to the
# PCBShield BOM workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Link-column fixes (rows 29-31: 8DIP / 4DIP sockets, screw terminals) ---
# Row 30 (4DIP Socket) was pointing at the wrong (ED08DT) DigiKey listing;
# point it at the same ED16DT listing used by row 29.
$ws.Range("G30").Value = "https://www.digikey.com/product-detail/en/on-shore-technology-inc/ED16DT/ED3046-5-ND/4147596"

# Row 31 (2 Pin Screw Terminals): swap the primary Link to the new DigiKey
# part, and keep the old SparkFun link around in the Resources column.
$ws.Range("G31").Value = "https://www.digikey.com/product-detail/en/on-shore-technology-inc/OSTVN02A150/ED10561-ND/1588862"
$ws.Range("H31").Value = "https://www.sparkfun.com/products/10571"

# --- Final quantities/costs: Printed Circuit Board row ---
# Replace the computed "Actual Cost" formula with the final, actual value.
$ws.Range("C32").Value = 20.47

# --- Un-hide & resize the helper columns (E: Total, F: New Cost) ---
$ws.Columns("E").Hidden = $false
$ws.Columns("F").Hidden = $false
$ws.Columns("F").ColumnWidth = 8.109375

# --- View state: drop the scrolled/zoomed-out review view, back to normal ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("G2").Select()

$wb.Save()
